# "first red box done" - fill in the worked example on the "Protect" sheet
# (VectorA = B - A style block) mirroring the already-completed "Hide" sheet,
# and make "Protect" the active tab/sheet.

$wb = $excel.ActiveWorkbook
$hide    = $wb.Worksheets.Item("Hide")
$protect = $wb.Worksheets.Item("Protect")

# --- Vector / magnitude / normalised / new-point block (rows 14-18) ---------
# Values + formulas mirror the Hide sheet's rows 10-14 (offset of +4 rows).
$protect.Range("J14").Value = "ab"
$protect.Range("K14").Formula = "=K13-K12"
$protect.Range("L14").Formula = "=L13-L12"

$protect.Range("J15").Value = "|ab|"
$protect.Range("K15").Formula = "=SQRT(K14*K14+L14*L14)"

$protect.Range("J16").Value = "u"
$protect.Range("K16").Formula = "=K14/K15"
$protect.Range("L16").Formula = "=L14/K15"

$protect.Range("K18").Formula = "=K12+K16*K17"
$protect.Range("L18").Formula = "=L12+L16*K17"

# Copy the number/cell formatting from the matching cells on the Hide sheet
# so the new cells pick up the same styling used there.
$hide.Range("J10:L10").Copy()
$protect.Range("J14:L14").PasteSpecial(-4122)
$hide.Range("J11:L11").Copy()
$protect.Range("J15:L15").PasteSpecial(-4122)
$hide.Range("J12:L12").Copy()
$protect.Range("J16:L16").PasteSpecial(-4122)

# M16 becomes a plain empty cell styled like the "Robot" label cell (M15).
$protect.Range("M15").Copy()
$protect.Range("M16").PasteSpecial(-4122)

# The placeholder cells that used to sit under the "Robot" mini-table are no
# longer needed now that the "VectorA"/"VectorR" box occupies that space.
$protect.Range("N16:P16").Clear()

# --- New "VectorA" / "VectorR" / "Sum product" worked box (rows 17-19) -----
$protect.Range("M17").Value = "VectorA"
$protect.Range("N17").Value = "A"
$protect.Range("O17").Value = 2
$protect.Range("P17").Value = 2

$protect.Range("M18").Value = "VectorR"
$protect.Range("N18").Value = "R"
$protect.Range("O18").Value = 5
$protect.Range("P18").Value = 3

$protect.Range("M19").Value = "Sum product"
$protect.Range("N19").Formula = "=SUMPRODUCT(O17:P17, O18:P18)"

# Style the new labels/values like the existing "A"/"B" lookup table above it.
$protect.Range("M15").Copy()
$protect.Range("M17:M19").PasteSpecial(-4122)

$protect.Range("O12").Copy()
$protect.Range("N17:P17").PasteSpecial(-4122)
$protect.Range("N18:P18").PasteSpecial(-4122)
$protect.Range("N19").PasteSpecial(-4122)

# Highlight the as-yet-uncalculated sum-product result cells in red ("Bad").
$protect.Range("O19:P19").Style = "Bad"

# --- Sheet view / active tab -------------------------------------------------
$hide.Select()
$hide.Range("J10:L12").Select()
$protect.Select()
$protect.Range("N23").Select()
